$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new bilingual test rows appended below the existing data (rows 1-15),
# written in the same order the shared-string table records them so the
# resulting workbook lines up with the authored fixture:
#   43 -> "这条可以被提取 + 1"
#   44 -> "this can be selected +2."
#   45 -> "这条可以被提取 + 2"
#   46 -> "this can be selected +1."
$ws.Range("A16").Value = "这条可以被提取 + 1"
$ws.Range("B17").Value = "this can be selected +2."
$ws.Range("A17").Value = "这条可以被提取 + 2"
$ws.Range("B16").Value = "this can be selected +1."
$ws.Range("C16").Value = "这条可以被提取 + 1"
$ws.Range("C17").Value = "这条可以被提取 + 2"

# Selection moves to G23, matching the final cursor position recorded in the
# commit's sheet1.xml.
$ws.Range("G23").Select() | Out-Null
